$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all data rows (2..111)
# from 2023-10-06 (45205) to 2023-10-07 (45206).
$ws.Range("C2:C111").Value = 45206
